# ES101.xlsx edit script
# 1) Insert a new field row ("schPayMed" / Payment Method / ComboBox) into the
#    "View Dfn" sheet's schFg field group, right before the existing
#    "bttDownload" row (i.e. new row 35, pushing everything from the old
#    row 35 onward down by one row).
# 2) Bump the field-group column count (L25) from 3 to 4.
# 3) Bump the revision/build timestamp in C14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("View Dfn")

# --- Insert the new row, shifting rows 35:55 down to 36:56 -----------------
$ws.Rows("35:35").Insert()

# The freshly inserted row doesn't inherit the surrounding "data" style, so
# copy number-format/fill/border/alignment from the row above (row 34, which
# carries the same style used throughout this field table) onto the new row.
$ws.Range("C34:P34").Copy()
$ws.Range("C35:P35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate the new field-definition row ----------------------------------
$ws.Range("D35").Value = "schPayMed"
$ws.Range("E35").Value = "Payment Method"
$ws.Range("K35").Value = "ComboBox"
$ws.Range("L35").Value = "dataUrl: getPayMed`r`nvalues: {""value"": ""id"", ""display"": ""tp""}"

# --- Update the field-group's field count (schFg now has 4 fields) ---------
$ws.Range("L25").Value = 4

# --- Bump the revision timestamp --------------------------------------------
$ws.Range("C14").Value = "20250603170900"
